$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 corresponds to the QLD entry in the "Days since Health Leadership
# faced the public on COVID" table. Update the "Last Date" and "News Link"
# text to reflect the newer article.
$ws.Range("B7").Value = 45122
$ws.Range("C7").Value = "https://www.abc.net.au/news/2023-07-15/influenza-a-b-c-d-explained-flu-season/102599454"
